$d = $word.ActiveDocument
$q1 = [char]0x2018   # U+2018 LEFT SINGLE QUOTATION MARK
$q2 = [char]0x2019   # U+2019 RIGHT SINGLE QUOTATION MARK

# ---------------------------------------------------------------------
# Change 1 ("Relay Control: ..." bullet):
#   "Relay Control: These signals control go to the Economizer board
#   to control the relays."
#   -> "Relay Control: These signals go to the Economizer board to
#   control the relays."
# i.e. delete the duplicated word "control " that precedes "go". Word
# splits the paragraph's single run around the edit point and drops
# its "_GoBack" bookmark there, plus an extra (pre-existing) run
# boundary right before "Relay".
# ---------------------------------------------------------------------

$dupRange = $d.Content
$found = $dupRange.Find.Execute("signals control go")
if (-not $found) { throw "Could not find 'signals control go' to fix the duplicated word" }

$delStart = $dupRange.Start + "signals ".Length
$delEnd = $delStart + "control ".Length
$delRange = $d.Range($delStart, $delEnd)
if ($delRange.Text -ne "control ") { throw "Unexpected text at deletion point: '$($delRange.Text)'" }
$delRange.Delete()

# Word keeps a single, document-wide "_GoBack" bookmark marking the
# most recent edit location; (re)adding it here both drops it from
# wherever it used to be and places it at this edit point.
$gobackRange = $d.Range($delStart, $delStart)
$d.Bookmarks.Add("_GoBack", $gobackRange)

# Split "Relay" into its own run by toggling a character property on
# and back off right on that word - this splits the run without
# altering the resulting formatting, matching the target layout.
$relayRange = $d.Content
$found = $relayRange.Find.Execute("Relay Control: These signals")
if (-not $found) { throw "Could not find the 'Relay Control' sentence" }
$relayStart = $relayRange.Start
$relayEnd = $relayStart + "Relay".Length
$relaySub = $d.Range($relayStart, $relayEnd)
$relaySub.Font.Bold = 1
$relaySub.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 2 (default-password bullet):
#   " 'adm" | (old "_GoBack" bookmark) | "in' user is 'password'. ..."
#   -> " 'admin' user is 'password'. ..." as a single run.
# Since "_GoBack" only ever exists once, adding it above already
# removed it from here; only the left-over run split around the old
# bookmark remains, so re-merge that text into a single run.
# ---------------------------------------------------------------------

$adminRange = $d.Content
$found = $adminRange.Find.Execute($q1 + "admin" + $q2)
if (-not $found) { throw "Could not find the '" + $q1 + "admin" + $q2 + "' text" }
$adminRange.Delete()
$adminRange.InsertBefore($q1 + "admin" + $q2)

Write-Output "done"
